# "question summary&proper nouns.xlsx" update - 12/10/2018 question list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("问题总结")

# Extend the existing question text in C12 with an extra clarifying sentence
$ws.Range("C12").Value = "螺纹钢和钢坯价格对运价影响有多大? 尤其是广州的螺纹钢和钢坯价格"

# That row now wraps onto two lines like the other question rows, so bump its height to match (34)
$ws.Rows.Item(12).RowHeight = 34

# Add two new question rows following the existing pattern used by rows 6/8/10/12:
#   col A -> date 12/10/2018 (serial 43444), col B -> asker "Le ", col C -> the question text
$ws.Range("A14").Value = 43444
$ws.Range("B14").Value = "Le "
$ws.Range("C14").Value = "动力煤是指的货物为动力煤，还是作为船的燃料？"
$ws.Rows.Item(14).RowHeight = 34

$ws.Range("A16").Value = 43444
$ws.Range("B16").Value = "Le "
$ws.Range("C16").Value = "铁矿石的港口库存情况（bloomberg上数据） 是否能对船价产生比较大的影响？"
$ws.Rows.Item(16).RowHeight = 34

# Match styling (fonts/alignment/number formats) of the row above for the two new rows
$ws.Range("A12:C12").Copy() | Out-Null
$ws.Range("A14:C14").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-apply values/height in case PasteSpecial touched them
$ws.Range("A14").Value = 43444
$ws.Range("B14").Value = "Le "
$ws.Range("C14").Value = "动力煤是指的货物为动力煤，还是作为船的燃料？"
$ws.Rows.Item(14).RowHeight = 34

$ws.Range("A16").Value = 43444
$ws.Range("B16").Value = "Le "
$ws.Range("C16").Value = "铁矿石的港口库存情况（bloomberg上数据） 是否能对船价产生比较大的影响？"
$ws.Rows.Item(16).RowHeight = 34

# Move the active selection to D6, matching the saved view state
$ws.Range("D6").Select() | Out-Null
